# Append newly-collected 45C thermal-curve size data (plate1, date 20250820)
# for wells A01-A12 and B01-B12, matching the USDA data addition described
# in the commit message, and update the sheet view scaling/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wells = @(
    "A01","A02","A03","A04","A05","A06","A07","A08","A09","A10","A11","A12",
    "B01","B02","B03","B04","B05","B06","B07","B08","B09","B10","B11","B12"
)

$lengths = @(
    7.626, 14.435, 5.887, 9.316, 8.84, 11.705, 7.209, 11.061, 7.479, 14.865, 12.222, 6.081,
    6.245, 10.406, 8.777, 9.045, 10.041, 9.992, 9.106, 4.716, 9.141, 7.043, 9.889, 8.456
)

$startRow = 122

for ($i = 0; $i -lt $wells.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = 20250820
    $ws.Cells.Item($row, 2).Value = 45
    $ws.Cells.Item($row, 3).Value = "plate1"
    $ws.Cells.Item($row, 4).Value = $wells[$i]
    $ws.Cells.Item($row, 5).Value = $lengths[$i]
}

# Update the view to reflect the new data extent/selection
$ws.Application.ActiveWindow.ScrollRow = 115
$ws.Range("C148").Select()
